# Enhance the "Audio Files Metadata" template:
#   - set explicit column widths on the data sheet
#   - store the customer_satisfaction / handle_time example values as text
#   - add a new "Instructions" sheet documenting every column

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Audio Files Metadata")

# --- Column widths for the "Audio Files Metadata" sheet ---
$ws1.Columns.Item(1).ColumnWidth = 25
$ws1.Range("B1:I1").EntireColumn.ColumnWidth = 18

# --- H2/I2/H3/I3 were stored as numbers; re-enter them as text ---
$ws1.Range("H2:I3").NumberFormat = "@"
$ws1.Range("H2").Value = "4"
$ws1.Range("I2").Value = "180"
$ws1.Range("H3").Value = "5"
$ws1.Range("I3").Value = "240"

# --- Add a new "Instructions" sheet right after "Audio Files Metadata" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Instructions"

$ws2.Columns.Item(1).ColumnWidth = 25
$ws2.Columns.Item(2).ColumnWidth = 70

$ws2.Range("A1").Value = "Audio File Metadata - Instructions"
$ws2.Range("A2").Formula = '=""'
$ws2.Range("A3").Value = "This template is used to provide metadata for audio files in batch uploads."
$ws2.Range("A4").Formula = '=""'
$ws2.Range("A5").Value = "Field Instructions:"

$ws2.Range("A6").Value = "filename"
$ws2.Range("B6").Value = "Must match the exact filename of the uploaded audio file (including extension)"

$ws2.Range("A7").Value = "language"
$ws2.Range("B7").Value = "Use one of: english, spanish, french, hindi, other"

$ws2.Range("A8").Value = "version"
$ws2.Range("B8").Value = "Version number or identifier of the call script/process used"

$ws2.Range("A9").Value = "call_date"
$ws2.Range("B9").Value = "Date of the call in YYYY-MM-DD format"

$ws2.Range("A10").Value = "call_type"
$ws2.Range("B10").Value = "Type of call, e.g., inbound, outbound, service, sales, etc."

$ws2.Range("A11").Value = "agent_id"
$ws2.Range("B11").Value = "ID of the agent who handled the call"

$ws2.Range("A12").Value = "call_id"
$ws2.Range("B12").Value = "Unique identifier of the call (if available)"

$ws2.Range("A13").Value = "customer_satisfaction"
$ws2.Range("B13").Value = "Customer satisfaction score, typically 1-5"

$ws2.Range("A14").Value = "handle_time"
$ws2.Range("B14").Value = "Call duration in seconds"

# Keep the primary data sheet as the active tab, same as before the edit
$ws1.Activate()

Write-Output "Done"
